# Mapeo de descarga de eprepago
# Update the data-driven test sheet "Datos" with the new e-prepago values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Row 2: usuario (D2) changes from "autotest32" to "autotest25"
$ws.Range("D2").Value = "autotest25"

# Row 2: numeroCuenta (O2) changes from "406-733040-20" to "406-714500-19"
$ws.Range("O2").Value = "406-714500-19"

# Row 2: valorDescarga (M2) changes from "100000" to "1000"
$ws.Range("M2").Value = "1000"

# Row 3: ID (A3) changes from 1 to 2
$ws.Range("A3").Value = 2

# Move the active selection off the table, onto the next empty row
$ws.Range("A4").Select()
